$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their existing text storage (avoid Excel
# auto-converting numeric-looking strings like "236.68" into numbers).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "90.120.71"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "3.109.47"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "236.68"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").Value = "617.68"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "1.09"
$ws.Range("E7").Value = "  +5.05%  "
$ws.Range("D8").Value = "0.367"
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.736"
$ws.Range("E10").Value = "  +3.03%  "
$ws.Range("B11").Value = "LidoStakedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D11").Value = "2.391.21"
$ws.Range("E11").Value = "  -22.52%  "
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "34.90"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "5.49"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("D16").Value = "90.183.93"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "3.698.89"
$ws.Range("D18").Value = "3.148.10"
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("D19").Value = "3.64"
$ws.Range("E19").Value = "  -4.28%  "
$ws.Range("D20").Value = "14.91"
$ws.Range("E20").Value = "  +8.00%  "
$ws.Range("D21").Value = "5.81"
$ws.Range("E21").Value = "  +6.71%  "
$ws.Range("D22").Value = "0.0000201"
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("D23").Value = "437.84"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "8.96"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").Value = "5.70"
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("D26").Value = "11.72"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "81.62"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").Value = "3.302.22"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "0.124"
$ws.Range("E30").Value = "  +43.39%  "
$ws.Range("D31").Value = "0.228"
$ws.Range("E31").Value = "  +19.22%  "
$ws.Range("E32").Value = "  +7.67%  "
$ws.Range("D33").Value = "9.21"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "0.169"
$ws.Range("E34").Value = "  +12.81%  "
$ws.Range("D35").Value = "0.923"
$ws.Range("E35").Value = "  -8.01%  "
$ws.Range("D36").Value = "7.59"
$ws.Range("E36").Value = "  +7.66%  "
$ws.Range("D37").Value = "25.96"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "499.19"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "1.92"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").Value = "1.34"
$ws.Range("E40").Value = "  +6.13%  "
$ws.Range("D41").Value = "0.444"
$ws.Range("E41").Value = "  +11.51%  "
$ws.Range("D42").Value = "3.71"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("D43").Value = "3.42"
$ws.Range("E43").Value = "  -7.66%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D46").Value = "158.43"
$ws.Range("E46").Value = "  +5.14%  "
$ws.Range("D47").Value = "0.712"
$ws.Range("E47").Value = "  +5.23%  "
$ws.Range("D48").Value = "1.90"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").Value = "  +3.82%  "
$ws.Range("D50").Value = "43.88"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "4.38"
$ws.Range("E51").Value = "  +0.72%  "
